# Update Name of Algo
# Apply the cell value changes as described by the diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3").Value = 6.785000000000001
$ws.Range("B4").Value = 6.669
$ws.Range("D6").Value = -7.509
$ws.Range("B7").Value = 7.181
$ws.Range("D7").Value = -7.569
$ws.Range("B8").Value = 6.867999999999999
$ws.Range("D8").Value = -7.419000000000001
$ws.Range("A11").Value = -21.952
$ws.Range("E11").Value = 12.704
$ws.Range("A12").Value = -21.631
$ws.Range("B12").Value = 6.381
$ws.Range("B14").Value = 7.637000000000002
$ws.Range("E14").Value = 12.832
$ws.Range("A15").Value = -21.018
$ws.Range("D19").Value = -8.242000000000001
$ws.Range("E19").Value = 12.526
$ws.Range("D21").Value = -7.222
$ws.Range("E21").Value = 12.938
$ws.Range("B22").Value = 6.929
$ws.Range("D24").Value = -7.532000000000001
$ws.Range("D25").Value = -7.860999999999999
